# ProductBurndown.xlsx — "reviews and reviews accept/reject"
#
# Sprint rows 23-27 get updated Actual-Hours / Planned-Hours figures
# (reviews took extra hours; some sprints' planned hours were
# re-estimated after accept/reject review), and a brand-new sprint row
# (41744 = 2014-04-14) is inserted before the summary row so the
# burndown keeps going.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1 - Product Burndown")

# ---------------------------------------------------------------
# 1. Insert a new row at 34 for the new sprint, pushing the old
#    summary row (SUM/AVERAGE) down to row 35. Copy the formatting
#    from row 33 (the last sprint row) onto the newly inserted row.
# ---------------------------------------------------------------
$ws.Range("A34:K34").EntireRow.Insert(-4121, 0)
$ws.Range("A33:K33").Copy()
$ws.Range("A34:K34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 2. Data edits to existing sprint rows (reviews / accept-reject
#    changed the actual hours booked and a couple of planned-hours
#    re-estimates).
# ---------------------------------------------------------------
$ws.Range("E23").Value = 12

$ws.Range("B24").Value = 12
$ws.Range("E24").Value = 3

$ws.Range("B25").Value = 0
$ws.Range("E25").Value = 0

$ws.Range("E26").Value = 0

$ws.Range("B27").Value = 6

$ws.Range("B31").Value = 0

$ws.Range("B32").Value = 0

# ---------------------------------------------------------------
# 3. Row 33 keeps its same formulas, written out explicitly now
#    that it is the last row of the (shrunken) shared-formula
#    blocks.
# ---------------------------------------------------------------
$ws.Range("F33").Formula = "=J32+C33"
$ws.Range("G33").Formula = "=H32"
$ws.Range("I33").Formula = "=I32+E32"
$ws.Range("J33").Formula = "=MAX(IF(OR(ISBLANK(D33),ISBLANK(E33)),F33-K32*B33,F33-D33),0)"
$ws.Range("K33").Formula = "=IF(OR(ISBLANK(D33),ISBLANK(E33)),K32,H33/(I33+E33))"

# ---------------------------------------------------------------
# 4. New sprint row 34 (2014-04-14, A34 serial 41744) with its own
#    data/formulas, following the same pattern as row 33.
# ---------------------------------------------------------------
$ws.Range("A34").Value = 41744
$ws.Range("B34").Value = 8
$ws.Range("C34").ClearContents()
$ws.Range("D34").ClearContents()
$ws.Range("E34").ClearContents()

$ws.Range("F34").Formula = "=J33+C34"
$ws.Range("G34").Formula = "=H33"
$ws.Range("H34").Formula = "=G34+D34"
$ws.Range("I34").Formula = "=I33+E33"
$ws.Range("J34").Formula = "=MAX(IF(OR(ISBLANK(D34),ISBLANK(E34)),F34-K33*B34,F34-D34),0)"
$ws.Range("K34").Formula = "=IF(OR(ISBLANK(D34),ISBLANK(E34)),K33,H34/(I34+E34))"

# ---------------------------------------------------------------
# 5. Row 48 (first row below the data block) loses its taller
#    20.65pt height now that the new trailing blank row 54 carries
#    that instead.
# ---------------------------------------------------------------
$ws.Range("A48").RowHeight = 20.45
$ws.Range("A54").RowHeight = 20.65

# ---------------------------------------------------------------
# 6. Recalculate and restore the selection the author ended on.
# ---------------------------------------------------------------
$wb.Application.Calculate()

$ws.Activate()
$ws.Range("L34").Select()
